# Swap match-data rows that were recorded in reversed order.
# Each pair below represents two fixtures played on the same date whose
# rows (columns B:AD - i.e. everything except the running index in column A)
# need to be swapped with one another.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(2, 3),
    @(34, 35),
    @(60, 61),
    @(64, 65),
    @(88, 89),
    @(122, 123),
    @(144, 145),
    @(162, 163),
    @(188, 189),
    @(212, 213),
    @(224, 225),
    @(228, 229)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Columns B (2) through AD (30) inclusive.
    $range1 = $ws.Range($ws.Cells.Item($r1, 2), $ws.Cells.Item($r1, 30))
    $range2 = $ws.Range($ws.Cells.Item($r2, 2), $ws.Cells.Item($r2, 30))

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
